# Updates the cryptos worksheet cell values to match the refreshed market data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-TextValue "D2" "64.380.47"
Set-TextValue "E2" "  +1.36%  "

Set-TextValue "D3" "2.625.92"
Set-TextValue "E3" "  +0.39%  "

Set-TextValue "E4" "  +0.08%  "

Set-TextValue "D5" "593.31"
Set-TextValue "E5" "  -0.26%  "

Set-TextValue "D6" "152.14"
Set-TextValue "E6" "  +1.32%  "

Set-TextValue "E7" "  +0.06%  "

Set-TextValue "E8" "  +0.04%  "

Set-TextValue "E9" "  +4.69%  "

Set-TextValue "E10" "  +3.63%  "

Set-TextValue "D11" "5.79"
Set-TextValue "E11" "  +2.03%  "

Set-TextValue "E12" "  +1.25%  "

Set-TextValue "D13" "28.47"
Set-TextValue "E13" "  +3.04%  "

Set-TextValue "D14" "3.100.72"
Set-TextValue "E14" "  +0.50%  "

Set-TextValue "B15" "WrappedBTC"
Set-TextValue "C15" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D15" "64.358.33"
Set-TextValue "E15" "  +1.58%  "

Set-TextValue "B16" "ShibaInu"
Set-TextValue "C16" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D16" "0.0000171"
Set-TextValue "E16" "  +13.09%  "

Set-TextValue "D17" "2.604.57"
Set-TextValue "E17" "  -0.41%  "

Set-TextValue "D18" "12.25"
Set-TextValue "E18" "  -0.65%  "

Set-TextValue "E19" "  +2.19%  "

Set-TextValue "D20" "350.41"
Set-TextValue "E20" "  +1.15%  "

Set-TextValue "D21" "7.13"
Set-TextValue "E21" "  +4.27%  "

Set-TextValue "D22" "0.999"
Set-TextValue "E22" "  +0.19%  "

Set-TextValue "D23" "67.52"
Set-TextValue "E23" "  +1.61%  "

Set-TextValue "D24" "1.69"
Set-TextValue "E24" "  -1.86%  "

Set-TextValue "D25" "9.30"
Set-TextValue "E25" "  +1.35%  "

Set-TextValue "E26" "  -1.19%  "

Set-TextValue "D27" "8.23"
Set-TextValue "E27" "  +1.13%  "

Set-TextValue "D28" "0.164"
Set-TextValue "E28" "  +1.44%  "

Set-TextValue "D29" "541.88"
Set-TextValue "E29" "  -1.88%  "

Set-TextValue "D30" "0.997"
Set-TextValue "E30" "  -0.23%  "

Set-TextValue "D31" "0.0₃0910"
Set-TextValue "E31" "  +7.88%  "

Set-TextValue "E32" "  +1.08%  "

Set-TextValue "D33" "1.81"
Set-TextValue "E33" "  +3.62%  "

Set-TextValue "D34" "5.62"
Set-TextValue "E34" "  +7.51%  "

Set-TextValue "E35" "  +0.78%  "

Set-TextValue "D36" "0.423"
Set-TextValue "E36" "  +2.23%  "

Set-TextValue "D37" "163.97"
Set-TextValue "E37" "  -2.35%  "

Set-TextValue "D38" "20.08"
Set-TextValue "E38" "  +3.46%  "

Set-TextValue "D39" "2.00"
Set-TextValue "E39" "  +3.72%  "

Set-TextValue "E40" "  +0.06%  "

Set-TextValue "E41" "  -0.04%  "

Set-TextValue "D42" "169.58"
Set-TextValue "E42" "  +1.96%  "

Set-TextValue "D43" "41.43"
Set-TextValue "E43" "  +4.32%  "

Set-TextValue "E44" "  +4.71%  "

Set-TextValue "D45" "23.24"
Set-TextValue "E45" "  +8.15%  "

Set-TextValue "D46" "0.0594"
Set-TextValue "E46" "  +1.29%  "

Set-TextValue "D47" "2.21"
Set-TextValue "E47" "  +11.19%  "

Set-TextValue "D48" "0.641"
Set-TextValue "E48" "  +2.03%  "

Set-TextValue "D49" "0.0251"
Set-TextValue "E49" "  +0.69%  "

Set-TextValue "D50" "0.0980"
Set-TextValue "E50" "  +1.45%  "

Set-TextValue "D51" "19.30"
Set-TextValue "E51" "  +0.38%  "
